$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: split the run containing [start, start+len) into its own run with
# no visible formatting change, by toggling Bold on then off. This forces
# the engine to materialize a distinct <w:r> for that span instead of
# re-coalescing it into its neighbour.
# ---------------------------------------------------------------------------
function Split-Run([int]$start, [int]$len) {
    $rg = $d.Range($start, $start + $len)
    $rg.Bold = 1
    $rg.Bold = 0
}

# ===========================================================================
# Hunk 1: "Place the salesfile.py and the salesfilerun.bat files in your
# Python Scripts folder"
#   -> "Place the salesfile" + "VX.X" + ".py and the salesfilerun" + "VX.X"
#      + ".bat files in your Python Scripts folder"
# ===========================================================================
$p1 = $d.Content
$p1.Find.Execute("Place the salesfile.py and the salesfilerun.bat files in your Python Scripts folder") | Out-Null
$p1Start = $p1.Start

$split1Offset = "Place the salesfile".Length
$split2Offset = "Place the salesfile.py and the salesfilerun".Length

# Insert the right-most marker first so the left-hand offset stays valid.
$ins2 = $d.Range($p1Start + $split2Offset, $p1Start + $split2Offset)
$ins2.InsertAfter("VX.X") | Out-Null

$ins1 = $d.Range($p1Start + $split1Offset, $p1Start + $split1Offset)
$ins1.InsertAfter("VX.X") | Out-Null

# Now carve each freshly-inserted "VX.X" into its own run.
Split-Run ($p1Start + $split1Offset) 4
Split-Run ($p1Start + $split2Offset + 4) 4

# ===========================================================================
# Hunk 2: "Double click on salesfilerun.bat to run the script and answer
# command line prompts" + <bookmark _GoBack/> + ". Enjoy!"
#   -> "Double click on salesfilerun" + "VX.X" + <bookmark _GoBack/>
#      + ".bat to run the script and answer command line prompts. Enjoy!"
# ===========================================================================
$p2 = $d.Content
$p2.Find.Execute("Double click on salesfilerun.bat to run the script and answer command line prompts") | Out-Null
$p2Start = $p2.Start

$split3Offset = "Double click on salesfilerun".Length

# The _GoBack bookmark currently sits right before ". Enjoy!" (i.e. right
# after "...command line prompts"). It needs to end up right before
# ".bat to run..." instead, i.e. immediately after our new "VX.X". Delete it
# and re-add it at the correct spot once the text has been inserted.
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

$ins3 = $d.Range($p2Start + $split3Offset, $p2Start + $split3Offset)
$ins3.InsertAfter("VX.X") | Out-Null

$bmPos = $d.Range($p2Start + $split3Offset + 4, $p2Start + $split3Offset + 4)
$d.Bookmarks.Add("_GoBack", $bmPos) | Out-Null

Split-Run ($p2Start + $split3Offset) 4

Write-Output "edit applied"
